$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B20").Value = 662
$ws.Range("C20").Value = "Maximum Width of Binary Tree"
$ws.Range("D20").Value = "Java"
$ws.Range("E20").Value = "Medium"

$ws.Range("B21").Value = "GFG"
$ws.Range("C21").Value = "Maximum Width of Binary Tree at any level."
$ws.Range("D21").Value = "java"
$ws.Range("E21").Value = "Easy"

$ws.Range("D26").Select()
